$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range('D2')
$c.NumberFormat = "@"
$c.Value = '64.284.56'
$c.Style = "Normal"
$ws.Range('E2').Value = '  +0.86%  '

$c = $ws.Range('D3')
$c.NumberFormat = "@"
$c.Value = '3.502.96'
$c.Style = "Normal"
$ws.Range('E3').Value = '  +0.21%  '

$ws.Range('E4').Value = '  -0.05%  '

$c = $ws.Range('D5')
$c.NumberFormat = "@"
$c.Value = '586.23'
$c.Style = "Normal"
$ws.Range('E5').Value = '  +0.47%  '

$c = $ws.Range('D6')
$c.NumberFormat = "@"
$c.Value = '134.40'
$c.Style = "Normal"
$ws.Range('E6').Value = '  +3.01%  '

$ws.Range('E7').Value = '  -0.06%  '

$ws.Range('E8').Value = '  +0.67%  '

$ws.Range('E9').Value = '  +1.35%  '

$ws.Range('E10').Value = '  +1.44%  '

$c = $ws.Range('D11')
$c.NumberFormat = "@"
$c.Value = '0.387'
$c.Style = "Normal"
$ws.Range('E11').Value = '  +2.05%  '

$c = $ws.Range('D12')
$c.NumberFormat = "@"
$c.Value = '4.101.31'
$c.Style = "Normal"
$ws.Range('E12').Value = '  +0.19%  '

$ws.Range('E13').Value = '  +3.82%  '

$ws.Range('E14').Value = '  +1.19%  '

$c = $ws.Range('D15')
$c.NumberFormat = "@"
$c.Value = '3.508.96'
$c.Style = "Normal"
$ws.Range('E15').Value = '  +1.39%  '

$c = $ws.Range('D16')
$c.NumberFormat = "@"
$c.Value = '25.90'
$c.Style = "Normal"
$ws.Range('E16').Value = '  -5.00%  '

$c = $ws.Range('D17')
$c.NumberFormat = "@"
$c.Value = '64.305.98'
$c.Style = "Normal"
$ws.Range('E17').Value = '  +0.70%  '

$c = $ws.Range('D18')
$c.NumberFormat = "@"
$c.Value = '9.91'
$c.Style = "Normal"
$ws.Range('E18').Value = '  +0.62%  '

$c = $ws.Range('D19')
$c.NumberFormat = "@"
$c.Value = '5.75'
$c.Style = "Normal"
$ws.Range('E19').Value = '  +2.47%  '

$c = $ws.Range('D20')
$c.NumberFormat = "@"
$c.Value = '13.68'
$c.Style = "Normal"
$ws.Range('E20').Value = '  -2.50%  '

$c = $ws.Range('D21')
$c.NumberFormat = "@"
$c.Value = '393.26'
$c.Style = "Normal"
$ws.Range('E21').Value = '  +2.47%  '

$c = $ws.Range('D22')
$c.NumberFormat = "@"
$c.Value = '0.572'
$c.Style = "Normal"
$ws.Range('E22').Value = '  -0.50%  '

$c = $ws.Range('D23')
$c.NumberFormat = "@"
$c.Value = '3.644.45'
$c.Style = "Normal"
$ws.Range('E23').Value = '  +0.22%  '

$c = $ws.Range('D24')
$c.NumberFormat = "@"
$c.Value = '74.24'
$c.Style = "Normal"
$ws.Range('E24').Value = '  +1.23%  '

$ws.Range('E25').Value = '  -0.07%  '

$ws.Range('E26').Value = '  -1.09%  '

$ws.Range('E27').Value = '  +1.00%  '

$ws.Range('E28').Value = '  +0.03%  '

$c = $ws.Range('D29')
$c.NumberFormat = "@"
$c.Value = '7.40'
$c.Style = "Normal"
$ws.Range('E29').Value = '  -0.90%  '

$ws.Range('E30').Value = '  -5.01%  '

$c = $ws.Range('D31')
$c.NumberFormat = "@"
$c.Value = '8.27'
$c.Style = "Normal"
$ws.Range('E31').Value = '  +0.31%  '

$ws.Range('E32').Value = '  +0.12%  '

$c = $ws.Range('D33')
$c.NumberFormat = "@"
$c.Value = '3.524.92'
$c.Style = "Normal"
$ws.Range('E33').Value = '  +0.56%  '

$ws.Range('E34').Value = '  +4.45%  '

$ws.Range('E35').Value = '  +0.01%  '

$ws.Range('E36').Value = '  +0.20%  '

$c = $ws.Range('D37')
$c.NumberFormat = "@"
$c.Value = '5.18'
$c.Style = "Normal"
$ws.Range('E37').Value = '  -2.40%  '

$c = $ws.Range('D38')
$c.NumberFormat = "@"
$c.Value = '1.56'
$c.Style = "Normal"
$ws.Range('E38').Value = '  +0.59%  '

$c = $ws.Range('D39')
$c.NumberFormat = "@"
$c.Value = '6.92'
$c.Style = "Normal"
$ws.Range('E39').Value = '  -0.14%  '

$c = $ws.Range('D40')
$c.NumberFormat = "@"
$c.Value = '163.81'
$c.Style = "Normal"
$ws.Range('E40').Value = '  +2.07%  '

$c = $ws.Range('D41')
$c.NumberFormat = "@"
$c.Value = '0.0783'
$c.Style = "Normal"
$ws.Range('E41').Value = '  -0.90%  '

$ws.Range('E42').Value = '  -0.46%  '

$ws.Range('E43').Value = '  -0.06%  '

$c = $ws.Range('D44')
$c.NumberFormat = "@"
$c.Value = '25.35'
$c.Style = "Normal"
$ws.Range('E44').Value = '  -4.25%  '

$ws.Range('E45').Value = '  +1.23%  '

$ws.Range('B46').Value = 'Stacks'
$ws.Range('C46').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$c = $ws.Range('D46')
$c.NumberFormat = "@"
$c.Value = '1.65'
$c.Style = "Normal"
$ws.Range('E46').Value = '  +2.20%  '

$ws.Range('B47').Value = 'ONDO'
$ws.Range('C47').Value = 'https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo'
$c = $ws.Range('D47')
$c.NumberFormat = "@"
$c.Value = '1.17'
$c.Style = "Normal"
$ws.Range('E47').Value = '  -2.85%  '

$c = $ws.Range('D48')
$c.NumberFormat = "@"
$c.Value = '2.464.57'
$c.Style = "Normal"
$ws.Range('E48').Value = '  +1.49%  '

$ws.Range('E49').Value = '  -0.44%  '

$c = $ws.Range('D50')
$c.NumberFormat = "@"
$c.Value = '0.898'
$c.Style = "Normal"
$ws.Range('E50').Value = '  +0.37%  '

$ws.Range('E51').Value = '  -1.56%  '

